$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recibo (D6) moves from the red/plain style to the green "proxy" style
# already used by row 7 (font color FF00B050 -> BGR long 0x50B000).
$ws.Range("D6").Font.Color = 5287936

# New "Faltan:" section below the diagram: Monto? / FormaPago? labels in
# bold, plus the "Faltan:" caption. Write the new shared strings in the
# same order they appear in the target sharedStrings table.
$ws.Range("B12").Value = "Monto?"
$ws.Range("B13").Value = "FormaPago?"
$ws.Range("A12").Value = "Faltan:"

$ws.Range("B12").Font.Bold = $true
$ws.Range("B13").Font.Bold = $true

# Leave the selection where the author ended up after typing the new rows.
[void]$ws.Range("C14").Select()
